# Saldo_guide.xlsx refresh: data was re-pulled a day later (IClientBalance
# export date rolled from 2024-09-11 to 2024-09-12), so every balance row's
# date column (G) advances by one day, and a handful of client balances
# (D/E/H) were updated with the newer figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every data row (2-274): bump the "date" column G from 45546 to 45547.
$ws.Range("G2:G274").Value = 45547

# Rows whose balances changed in the refreshed extract.
# H (Saldo Atual / current balance) = D (Debito) + E (Credito) in every row.
$ws.Range("D8").Value = 35189.58
$ws.Range("H8").Value = 46487.76

$ws.Range("E55").Value = 1108.72
$ws.Range("H55").Value = 1108.72

$ws.Range("E58").Value = 49.26
$ws.Range("H58").Value = 49.26

$ws.Range("E101").Value = 193.56
$ws.Range("H101").Value = 193.56

$ws.Range("E112").Value = 447.28
$ws.Range("H112").Value = 447.28

$ws.Range("E113").Value = 596.45000000000005
$ws.Range("H113").Value = 596.45000000000005

$ws.Range("E197").Value = 71.36
$ws.Range("H197").Value = 71.36

$ws.Range("E232").Value = 370.57
$ws.Range("H232").Value = 370.57

$ws.Range("E255").Value = 720.39
$ws.Range("H255").Value = 720.39

# Sheet / workbook metadata: the export file name embeds the pull timestamp.
$wb.Sheets.Item(1).Name = "IClientBalance-20240912-084643-"
